# Update countries & provincias Spain
# Applies the data refresh captured in the diff:
#  - Updated "Datos actualizados" timestamp string
#  - Updated case numbers for Alemania (row 12) and Israel (row 48)
#  - Re-sorted / updated adjacent country rows whose case totals crossed over:
#       Bolivia / Barein        (rows 50-51)
#       Tayikistan / Uzbekistan (rows 76-77)
#       Islas Virgenes Britanicas / Papua Nueva Guinea (rows 213-214)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Junio de 2020 a las 08:41"

# --- Simple numeric updates (no reordering) ---

# Row 12: Alemania
$ws.Range("D12").Value = 172200
$ws.Range("E12").Value = 6356

# Row 48: Israel
$ws.Range("B48").Value = 19008
$ws.Range("C48").Value = 36
$ws.Range("D48").Value = 15360
$ws.Range("E48").Value = 3348

# --- Row 50/51 swap: Barein now ranks above Bolivia ---
$ws.Range("A50").Value = "Barein"
$ws.Range("B50").Value = 18227
$ws.Range("C50").Value = 514
$ws.Range("D50").Value = 12818
$ws.Range("E50").Value = 5370
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 39

$ws.Range("A51").Value = "Bolivia"
$ws.Range("B51").Value = 17842
$ws.Range("C51").Value = 913
$ws.Range("D51").Value = 2768
$ws.Range("E51").Value = 14489
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 26
$ws.Range("H51").Value = 585

# --- Row 76/77 swap: Uzbekistan now ranks above Tayikistan ---
$ws.Range("A76").Value = "Uzbekistan"
$ws.Range("B76").Value = 4994
$ws.Range("C76").Value = 28
$ws.Range("D76").Value = 3874
$ws.Range("E76").Value = 1101
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 19

$ws.Range("A77").Value = "Tayikistan"
$ws.Range("B77").Value = 4971
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 3288
$ws.Range("E77").Value = 1633
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 50

# --- Row 213/214 swap: Papua Nueva Guinea now ranks above Islas Virgenes Britanicas ---
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 8
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 7
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
